$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.4654124536498146
$ws.Range("D2").Value = -0.1932875571502286
$ws.Range("E2").Value = -0.4647271069084276
$ws.Range("F2").Value = 0.4123650414601658
$ws.Range("G2").Value = 0.0318876795507182
$ws.Range("H2").Value = -0.06818235272941092
$ws.Range("I2").Value = 0.0146936267745071
$ws.Range("J2").Value = 0.02268172272689091
$ws.Range("K2").Value = -0.003558542234168937
$ws.Range("L2").Value = -0.09407567230268921
$ws.Range("M2").Value = -0.03695851583406334
$ws.Range("N2").Value = -0.1359832959331837
$ws.Range("O2").Value = 0.2270926523706095
$ws.Range("P2").Value = -0.003746606986427946
$ws.Range("Q2").Value = 0.07204099216396866
$ws.Range("R2").Value = -0.1668656594626378
$ws.Range("S2").Value = 0.05952158208632834
$ws.Range("T2").Value = 0.00811308845235381
$ws.Range("U2").Value = -0.6684836979347918
$ws.Range("V2").Value = -0.07785362341449366
$ws.Range("W2").Value = -0.007372349489397957
$ws.Range("X2").Value = 0.3748385873543494
$ws.Range("Y2").Value = 0.1327970431881728
$ws.Range("Z2").Value = -0.4172270769083076
$ws.Range("AA2").Value = -0.08634927339709358
$ws.Range("AB2").Value = 0.08748620594482379
$ws.Range("AC2").Value = 0.02072542690170761
$ws.Range("AD2").Value = 0.06508202032808132
$ws.Range("AE2").Value = -0.001210468841875368
$ws.Range("AF2").Value = -0.06792823971295885
$ws.Range("AG2").Value = 0.01548975795903184
$ws.Range("AH2").Value = -0.02966575066300265
$ws.Range("B3").Value = -0.4654124536498146
$ws.Range("D3").Value = -0.1445168660674643
$ws.Range("E3").Value = 0.2620742322969292
$ws.Range("F3").Value = 0.1749964119856479
$ws.Range("G3").Value = -0.0234313897255589
$ws.Range("H3").Value = -0.1649794439177757
$ws.Range("I3").Value = 0.2832150528602114
$ws.Range("J3").Value = -0.01040951363805455
$ws.Range("K3").Value = -0.01913326053304213
$ws.Range("L3").Value = 0.006965115860463442
$ws.Range("M3").Value = -0.01088317953271813
$ws.Range("N3").Value = -0.002950091800367202
$ws.Range("O3").Value = -0.3351573086292345
$ws.Range("P3").Value = 0.04471937887751551
$ws.Range("Q3").Value = 0.3354629738518954
$ws.Range("R3").Value = -0.1081049764199057
$ws.Range("S3").Value = -0.01022020088080352
$ws.Range("T3").Value = 0.004378961515846064
$ws.Range("U3").Value = 0.5880355681422725
$ws.Range("V3").Value = 0.0004262417049668199
$ws.Range("W3").Value = -0.08759631838527354
$ws.Range("X3").Value = -0.3724452017808071
$ws.Range("Y3").Value = -0.1413225492901972
$ws.Range("Z3").Value = -0.1562407209628839
$ws.Range("AA3").Value = -0.2096079104316417
$ws.Range("AB3").Value = 0.0098687434749739
$ws.Range("AC3").Value = 0.04431540926163705
$ws.Range("AD3").Value = -0.06873992295969183
$ws.Range("AE3").Value = 0.004956019824079296
$ws.Range("AF3").Value = 0.02240706562826251
$ws.Range("AG3").Value = -0.3702358489433958
$ws.Range("AH3").Value = -0.1697342469369877
$ws.Range("B4").Value = -0.1932875571502286
$ws.Range("C4").Value = -0.1445168660674643
$ws.Range("E4").Value = 0.033562694250777
$ws.Range("F4").Value = -0.1277623030492122
$ws.Range("G4").Value = -0.04655788223152892
$ws.Range("H4").Value = -0.01044570578282313
$ws.Range("I4").Value = -0.07472593090372362
$ws.Range("J4").Value = -0.0298815595262381
$ws.Range("K4").Value = -0.09481573926295706
$ws.Range("L4").Value = 0.1044145296581186
$ws.Range("M4").Value = 0.1883777135108541
$ws.Range("N4").Value = 0.7066350185400742
$ws.Range("O4").Value = 0.06336409345637382
$ws.Range("P4").Value = 0.01116417265669063
$ws.Range("Q4").Value = 0.05697517190068761
$ws.Range("R4").Value = 0.4551829407317629
$ws.Range("S4").Value = 0.3285715862863451
$ws.Range("T4").Value = -0.01776487105948424
$ws.Range("U4").Value = -0.05932113328453314
$ws.Range("V4").Value = -0.0337189828759315
$ws.Range("W4").Value = 0.09931642926571706
$ws.Range("X4").Value = 0.006573242292969171
$ws.Range("Y4").Value = -0.05067601070404282
$ws.Range("Z4").Value = 0.1350391321565286
$ws.Range("AA4").Value = 0.05924241296965188
$ws.Range("AB4").Value = 0.03140758163032652
$ws.Range("AC4").Value = 0.01867927471709887
$ws.Range("AD4").Value = -0.08909517238068952
$ws.Range("AE4").Value = 0.0447559550238201
$ws.Range("AF4").Value = 0.009172548690194762
$ws.Range("AG4").Value = 0.06278175512702051
$ws.Range("AH4").Value = 0.001792327169308677
$ws.Range("B5").Value = -0.4647271069084276
$ws.Range("C5").Value = 0.2620742322969292
$ws.Range("D5").Value = 0.033562694250777
$ws.Range("F5").Value = -0.4920840483361933
$ws.Range("G5").Value = 0.07481444325777303
$ws.Range("H5").Value = -0.1381125844503378
$ws.Range("I5").Value = 0.1595471661886648
$ws.Range("J5").Value = 0.006848859395437582
$ws.Range("K5").Value = 0.02666276265105061
$ws.Range("L5").Value = -0.03219113676454706
$ws.Range("M5").Value = -0.05417714470857883
$ws.Range("N5").Value = -0.01839204156816627
$ws.Range("O5").Value = -0.2279471517886072
$ws.Range("P5").Value = 0.008653762615050461
$ws.Range("Q5").Value = 0.1519687998751995
$ws.Range("R5").Value = -0.1604484177936712
$ws.Range("S5").Value = 0.1405433141732567
$ws.Range("T5").Value = -0.002990219960879844
$ws.Range("U5").Value = 0.3868244272977092
$ws.Range("V5").Value = 0.02311958047832191
$ws.Range("W5").Value = 0.05188033552134209
$ws.Range("X5").Value = -0.06261567446269785
$ws.Range("Y5").Value = -0.07953948615794464
$ws.Range("Z5").Value = 0.4772423889695559
$ws.Range("AA5").Value = 0.5720074400297601
$ws.Range("AB5").Value = 0.1603518414073656
$ws.Range("AC5").Value = 0.01144132576530306
$ws.Range("AD5").Value = -0.09183396733586935
$ws.Range("AE5").Value = -0.06449228996915987
$ws.Range("AF5").Value = 0.0216188064752259
$ws.Range("AG5").Value = -0.143810015240061
$ws.Range("AH5").Value = 0.03881410725642902
$ws.Range("B6").Value = 0.4123650414601658
$ws.Range("C6").Value = 0.1749964119856479
$ws.Range("D6").Value = -0.1277623030492122
$ws.Range("E6").Value = -0.4920840483361933
$ws.Range("G6").Value = -0.02597415189660759
$ws.Range("H6").Value = 0.03115366061464246
$ws.Range("I6").Value = 0.07042165768663075
$ws.Range("J6").Value = 0.04512200448801795
$ws.Range("K6").Value = -0.03351344605378422
$ws.Range("L6").Value = -0.03838239352957412
$ws.Range("M6").Value = -0.04839590558362233
$ws.Range("N6").Value = -0.002676778707114828
$ws.Range("O6").Value = 0.01550291001164005
$ws.Range("P6").Value = -0.01825399301597206
$ws.Range("Q6").Value = -0.02707604430417722
$ws.Range("R6").Value = 0.02484777939111757
$ws.Range("S6").Value = 0.003350797403189613
$ws.Range("T6").Value = 0.03735394141576566
$ws.Range("U6").Value = 0.0274940299761199
$ws.Range("V6").Value = -0.06279106716426866
$ws.Range("W6").Value = -0.07668990675962704
$ws.Range("X6").Value = 0.0374810459241837
$ws.Range("Y6").Value = -0.01559910239640959
$ws.Range("Z6").Value = -0.9642329449317797
$ws.Range("AA6").Value = -0.1017374629498518
$ws.Range("AB6").Value = -0.04999421597686391
$ws.Range("AC6").Value = 0.2516081744326977
$ws.Range("AD6").Value = -0.002114312457249829
$ws.Range("AE6").Value = -0.003836463345853383
$ws.Range("AF6").Value = 0.009624518498073992
$ws.Range("AG6").Value = -0.04217949671798688
$ws.Range("AH6").Value = -0.07514545258181032
$ws.Range("B7").Value = 0.0318876795507182
$ws.Range("C7").Value = -0.0234313897255589
$ws.Range("D7").Value = -0.04655788223152892
$ws.Range("E7").Value = 0.07481444325777303
$ws.Range("F7").Value = -0.02597415189660759
$ws.Range("B8").Value = -0.06818235272941092
$ws.Range("C8").Value = -0.1649794439177757
$ws.Range("D8").Value = -0.01044570578282313
$ws.Range("E8").Value = -0.1381125844503378
$ws.Range("F8").Value = 0.03115366061464246
$ws.Range("B9").Value = 0.0146936267745071
$ws.Range("C9").Value = 0.2832150528602114
$ws.Range("D9").Value = -0.07472593090372362
$ws.Range("E9").Value = 0.1595471661886648
$ws.Range("F9").Value = 0.07042165768663075
$ws.Range("B10").Value = 0.02268172272689091
$ws.Range("C10").Value = -0.01040951363805455
$ws.Range("D10").Value = -0.0298815595262381
$ws.Range("E10").Value = 0.006848859395437582
$ws.Range("F10").Value = 0.04512200448801795
$ws.Range("B11").Value = -0.003558542234168937
$ws.Range("C11").Value = -0.01913326053304213
$ws.Range("D11").Value = -0.09481573926295706
$ws.Range("E11").Value = 0.02666276265105061
$ws.Range("F11").Value = -0.03351344605378422
$ws.Range("B12").Value = -0.09407567230268921
$ws.Range("C12").Value = 0.006965115860463442
$ws.Range("D12").Value = 0.1044145296581186
$ws.Range("E12").Value = -0.03219113676454706
$ws.Range("F12").Value = -0.03838239352957412
$ws.Range("B13").Value = -0.03695851583406334
$ws.Range("C13").Value = -0.01088317953271813
$ws.Range("D13").Value = 0.1883777135108541
$ws.Range("E13").Value = -0.05417714470857883
$ws.Range("F13").Value = -0.04839590558362233
$ws.Range("B14").Value = -0.1359832959331837
$ws.Range("C14").Value = -0.002950091800367202
$ws.Range("D14").Value = 0.7066350185400742
$ws.Range("E14").Value = -0.01839204156816627
$ws.Range("F14").Value = -0.002676778707114828
$ws.Range("B15").Value = 0.2270926523706095
$ws.Range("C15").Value = -0.3351573086292345
$ws.Range("D15").Value = 0.06336409345637382
$ws.Range("E15").Value = -0.2279471517886072
$ws.Range("F15").Value = 0.01550291001164005
$ws.Range("B16").Value = -0.003746606986427946
$ws.Range("C16").Value = 0.04471937887751551
$ws.Range("D16").Value = 0.01116417265669063
$ws.Range("E16").Value = 0.008653762615050461
$ws.Range("F16").Value = -0.01825399301597206
$ws.Range("B17").Value = 0.07204099216396866
$ws.Range("C17").Value = 0.3354629738518954
$ws.Range("D17").Value = 0.05697517190068761
$ws.Range("E17").Value = 0.1519687998751995
$ws.Range("F17").Value = -0.02707604430417722
$ws.Range("B18").Value = -0.1668656594626378
$ws.Range("C18").Value = -0.1081049764199057
$ws.Range("D18").Value = 0.4551829407317629
$ws.Range("E18").Value = -0.1604484177936712
$ws.Range("F18").Value = 0.02484777939111757
$ws.Range("B19").Value = 0.05952158208632834
$ws.Range("C19").Value = -0.01022020088080352
$ws.Range("D19").Value = 0.3285715862863451
$ws.Range("E19").Value = 0.1405433141732567
$ws.Range("F19").Value = 0.003350797403189613
$ws.Range("B20").Value = 0.00811308845235381
$ws.Range("C20").Value = 0.004378961515846064
$ws.Range("D20").Value = -0.01776487105948424
$ws.Range("E20").Value = -0.002990219960879844
$ws.Range("F20").Value = 0.03735394141576566
$ws.Range("B21").Value = -0.6684836979347918
$ws.Range("C21").Value = 0.5880355681422725
$ws.Range("D21").Value = -0.05932113328453314
$ws.Range("E21").Value = 0.3868244272977092
$ws.Range("F21").Value = 0.0274940299761199
$ws.Range("B22").Value = -0.07785362341449366
$ws.Range("C22").Value = 0.0004262417049668199
$ws.Range("D22").Value = -0.0337189828759315
$ws.Range("E22").Value = 0.02311958047832191
$ws.Range("F22").Value = -0.06279106716426866
$ws.Range("B23").Value = -0.007372349489397957
$ws.Range("C23").Value = -0.08759631838527354
$ws.Range("D23").Value = 0.09931642926571706
$ws.Range("E23").Value = 0.05188033552134209
$ws.Range("F23").Value = -0.07668990675962704
$ws.Range("B24").Value = 0.3748385873543494
$ws.Range("C24").Value = -0.3724452017808071
$ws.Range("D24").Value = 0.006573242292969171
$ws.Range("E24").Value = -0.06261567446269785
$ws.Range("F24").Value = 0.0374810459241837
$ws.Range("B25").Value = 0.1327970431881728
$ws.Range("C25").Value = -0.1413225492901972
$ws.Range("D25").Value = -0.05067601070404282
$ws.Range("E25").Value = -0.07953948615794464
$ws.Range("F25").Value = -0.01559910239640959
$ws.Range("B26").Value = -0.4172270769083076
$ws.Range("C26").Value = -0.1562407209628839
$ws.Range("D26").Value = 0.1350391321565286
$ws.Range("E26").Value = 0.4772423889695559
$ws.Range("F26").Value = -0.9642329449317797
$ws.Range("B27").Value = -0.08634927339709358
$ws.Range("C27").Value = -0.2096079104316417
$ws.Range("D27").Value = 0.05924241296965188
$ws.Range("E27").Value = 0.5720074400297601
$ws.Range("F27").Value = -0.1017374629498518
$ws.Range("B28").Value = 0.08748620594482379
$ws.Range("C28").Value = 0.0098687434749739
$ws.Range("D28").Value = 0.03140758163032652
$ws.Range("E28").Value = 0.1603518414073656
$ws.Range("F28").Value = -0.04999421597686391
$ws.Range("B29").Value = 0.02072542690170761
$ws.Range("C29").Value = 0.04431540926163705
$ws.Range("D29").Value = 0.01867927471709887
$ws.Range("E29").Value = 0.01144132576530306
$ws.Range("F29").Value = 0.2516081744326977
$ws.Range("B30").Value = 0.06508202032808132
$ws.Range("C30").Value = -0.06873992295969183
$ws.Range("D30").Value = -0.08909517238068952
$ws.Range("E30").Value = -0.09183396733586935
$ws.Range("F30").Value = -0.002114312457249829
$ws.Range("B31").Value = -0.001210468841875368
$ws.Range("C31").Value = 0.004956019824079296
$ws.Range("D31").Value = 0.0447559550238201
$ws.Range("E31").Value = -0.06449228996915987
$ws.Range("F31").Value = -0.003836463345853383
$ws.Range("B32").Value = -0.06792823971295885
$ws.Range("C32").Value = 0.02240706562826251
$ws.Range("D32").Value = 0.009172548690194762
$ws.Range("E32").Value = 0.0216188064752259
$ws.Range("F32").Value = 0.009624518498073992
$ws.Range("B33").Value = 0.01548975795903184
$ws.Range("C33").Value = -0.3702358489433958
$ws.Range("D33").Value = 0.06278175512702051
$ws.Range("E33").Value = -0.143810015240061
$ws.Range("F33").Value = -0.04217949671798688
$ws.Range("B34").Value = -0.02966575066300265
$ws.Range("C34").Value = -0.1697342469369877
$ws.Range("D34").Value = 0.001792327169308677
$ws.Range("E34").Value = 0.03881410725642902
$ws.Range("F34").Value = -0.07514545258181032
